$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6490908265113831
$ws.Range("B1").Value = 2.059133529663086
$ws.Range("C1").Value = 2.640430688858032
$ws.Range("D1").Value = 0.6737529039382935
$ws.Range("E1").Value = 0.7557560801506042
